$d = $word.ActiveDocument

$replacements = @(
    @("80-16=64", "48-47=1"),
    @("98-3=95", "58+18=76"),
    @("71+3=74", "87-68=19"),
    @("25+12=37", "48+29=77"),
    @("90-63=27", "79+4=83"),
    @("37-22=15", "32-4=28"),
    @("27+10=37", "13+38=51"),
    @("36+44=80", "57+13=70"),
    @("36-7=29", "66-50=16"),
    @("9+14=23", "8+16=24"),
    @("99-44=55", "64+23=87"),
    @("35-34=1", "9+43=52"),
    @("58-6=52", "1+73=74"),
    @("8+41=49", "98-43=55"),
    @("69-27=42", "99-42=57"),
    @("62-54=8", "38+19=57"),
    @("95-4=91", "76+3=79"),
    @("90-19=71", "15+3=18"),
    @("11+46=57", "0+3=3"),
    @("53-1=52", "14-13=1"),
    @("92-65=27", "16+77=93"),
    @("2+33=35", "46+5=51"),
    @("66-34=32", "79-26=53"),
    @("98-19=79", "33+61=94"),
    @("20+29=49", "19-19=0"),
    @("12+81=93", "22-6=16"),
    @("75-10=65", "86-8=78"),
    @("23+45=68", "68-18=50"),
    @("42-2=40", "67-23=44"),
    @("24-8=16", "35-10=25"),
    @("7+4=11", "22+15=37"),
    @("83-15=68", "26+22=48"),
    @("1+34=35", "41-39=2"),
    @("78-69=9", "22+73=95"),
    @("79-40=39", "67-23=44"),
    @("81-78=3", "5+56=61"),
    @("47+0=47", "62-13=49"),
    @("88-67=21", "40-23=17"),
    @("24+69=93", "41-23=18"),
    @("92-12=80", "84-82=2"),
    @("64+27=91", "4+63=67"),
    @("28+5=33", "32-30=2"),
    @("1+29=30", "10-0=10"),
    @("50+37=87", "65-26=39"),
    @("20+73=93", "43-8=35"),
    @("77-43=34", "58+30=88"),
    @("26+70=96", "89-31=58"),
    @("79+0=79", "96-45=51"),
    @("19+45=64", "28+27=55"),
    @("38+18=56", "26+56=82"),
    @("59-31=28", "91-45=46"),
    @("25+22=47", "80+10=90"),
    @("12+14=26", "16+5=21"),
    @("78-37=41", "17+5=22"),
    @("36+1=37", "73-6=67"),
    @("45+48=93", "98-76=22"),
    @("53-23=30", "40+51=91"),
    @("0+42=42", "27+19=46"),
    @("95-57=38", "13+12=25"),
    @("18-13=5", "28+70=98"),
    @("43+53=96", "10+36=46"),
    @("63+28=91", "65+6=71"),
    @("95-75=20", "50-35=15"),
    @("14+6=20", "16-11=5"),
    @("59-48=11", "21+23=44"),
    @("11+15=26", "62-16=46"),
    @("66-19=47", "54+35=89"),
    @("0+37=37", "87-37=50"),
    @("81-8=73", "40-33=7"),
    @("6+65=71", "43-23=20"),
    @("38-0=38", "99-6=93"),
    @("41+33=74", "26+27=53"),
    @("79-55=24", "63+11=74"),
    @("42-4=38", "0+73=73"),
    @("67-6=61", "1+22=23"),
    @("77-45=32", "44-9=35"),
    @("62-37=25", "19+38=57"),
    @("71-12=59", "27+1=28"),
    @("56+41=97", "15+2=17"),
    @("76-16=60", "29+58=87"),
    @("54+25=79", "75-28=47"),
    @("50-48=2", "93-81=12"),
    @("49-34=15", "16+76=92"),
    @("73-28=45", "28+17=45"),
    @("7+21=28", "87-26=61"),
    @("63-21=42", "79-45=34"),
    @("0+38=38", "12+45=57"),
    @("21+34=55", "20+70=90"),
    @("62-24=38", "81-73=8"),
    @("97-87=10", "34+8=42"),
    @("15+43=58", "2+10=12"),
    @("0+43=43", "77+17=94"),
    @("36+61=97", "54-12=42"),
    @("83-54=29", "77-66=11"),
    @("28+20=48", "16+41=57"),
    @("90-84=6", "79-53=26"),
    @("65+7=72", "21+7=28"),
    @("34+59=93", "11+65=76"),
    @("37+9=46", "2+29=31"),
    @("5+39=44", "87-21=66"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
